# Butte recapture EDI: drop the "mort" column (J) from the "Recapture_EDI"
# sheet, refreshed from a new query. Every column from K onward shifts one
# position left (K->J, L->K, ... W->V) and the now-unused last column (W)
# is cleared. Only the first 16 rows (header + 15 data rows) of the new
# query have data, so rows 17-27 (the old extra data rows) are cleared out
# entirely. The workbook-level named range "Recapture_EDI" is updated to
# match the new, smaller extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recapture_EDI")

$lastDataRow = 16
$firstShiftCol = 11   # K
$lastShiftCol = 23    # W

for ($r = 1; $r -le $lastDataRow; $r++) {
    for ($c = $firstShiftCol; $c -le $lastShiftCol; $c++) {
        $srcVal = $ws.Cells.Item($r, $c).Value2
        $ws.Cells.Item($r, $c - 1).Value2 = $srcVal
    }
    $ws.Cells.Item($r, $lastShiftCol).Value2 = $null
}

$totalRows = 27
$totalCols = 23
for ($r = $lastDataRow + 1; $r -le $totalRows; $r++) {
    for ($c = 1; $c -le $totalCols; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $null
    }
}

$nm = $wb.Names.Item("Recapture_EDI")
$nm.RefersTo = "='Recapture_EDI'!`$A`$1:`$V`$16"
